$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2307.8125
$ws.Range("J112").Value = 2305.923
$ws.Range("L112").Value = 6917.768999999999
$ws.Range("N112").Value = -9133.769
$ws.Range("H116").Value = 5265.0347
$ws.Range("I116").Value = 5028.5264
$ws.Range("J116").Value = 5714.4
$ws.Range("K116").Value = 5028.5264
$ws.Range("L116").Value = 5714.4
$ws.Range("M116").Value = -1586.5264
$ws.Range("N116").Value = -12598.4
$ws.Range("H124").Value = 99999
$ws.Range("J124").Value = 99999
$ws.Range("L124").Value = 99999
$ws.Range("N124").Value = -109819
$ws.Range("H131").Value = 25790
$ws.Range("J131").Value = 49996.75
$ws.Range("L131").Value = 149990.25
$ws.Range("N131").Value = -160070.25
$ws.Range("H132").Value = 3340.652
$ws.Range("I132").Value = 3340.652
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10021.956
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7491.956
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 110000
$ws.Range("J133").Value = 110000
$ws.Range("L133").Value = 110000
$ws.Range("N133").Value = -120120
$ws.Range("H135").Value = 1829.875
$ws.Range("I135").Value = 1829.875
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 16468.875
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -13933.875
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 2232.7646
$ws.Range("J137").Value = 3255.7778
$ws.Range("L137").Value = 9767.3334
$ws.Range("N137").Value = -14867.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 10168
$ws.Range("I26").Value = 10168
$ws.Range("K26").Value = 10168
$ws.Range("M26").Value = -9838
$ws.Range("H32").Value = 33240.582
$ws.Range("I32").Value = 35567.688
$ws.Range("K32").Value = 35567.688
$ws.Range("M32").Value = -35280.688
$ws.Range("H61").Value = 2148.0908
$ws.Range("I61").Value = 2112.9
$ws.Range("K61").Value = 2112.9
$ws.Range("M61").Value = -1900.9
$ws.Range("H94").Value = 33333.332
$ws.Range("J94").Value = 33333.332
$ws.Range("L94").Value = 33333.332
$ws.Range("N94").Value = -35135.332
$ws.Range("H136").Value = 2148.0908
$ws.Range("I136").Value = 2112.9
$ws.Range("K136").Value = 6338.700000000001
$ws.Range("M136").Value = -3788.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 33833.332
$ws.Range("J6").Value = 33833.332
$ws.Range("L6").Value = 33833.332
$ws.Range("N6").Value = -34059.332
$ws.Range("H13").Value = 45000
$ws.Range("J13").Value = 45000
$ws.Range("L13").Value = 45000
$ws.Range("N13").Value = -45336
$ws.Range("H102").Value = 10411.444
$ws.Range("I102").Value = 10411.444
$ws.Range("K102").Value = 10411.444
$ws.Range("M102").Value = -7166.444
$ws.Range("H105").Value = 2800.4583
$ws.Range("I105").Value = 2942.7896
$ws.Range("K105").Value = 2942.7896
$ws.Range("M105").Value = -1195.7896
$ws.Range("H106").Value = 29999
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H122").Value = 93999.5
$ws.Range("J122").Value = 93999.5
$ws.Range("L122").Value = 93999.5
$ws.Range("N122").Value = -103799.5
$ws.Range("H124").Value = 45000
$ws.Range("J124").Value = 45000
$ws.Range("L124").Value = 45000
$ws.Range("N124").Value = -54820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1859.8462
$ws.Range("I22").Value = 338.5
$ws.Range("J22").Value = 3163.8572
$ws.Range("K22").Value = 338.5
$ws.Range("L22").Value = 3163.8572
$ws.Range("M22").Value = 11.5
$ws.Range("N22").Value = -3863.8572
$ws.Range("H31").Value = 2233.9656
$ws.Range("I31").Value = 1731.5714
$ws.Range("J31").Value = 3552.75
$ws.Range("K31").Value = 1731.5714
$ws.Range("L31").Value = 3552.75
$ws.Range("M31").Value = -1436.5714
$ws.Range("N31").Value = -4142.75
$ws.Range("H34").Value = 2233.9656
$ws.Range("I34").Value = 1731.5714
$ws.Range("J34").Value = 3552.75
$ws.Range("K34").Value = 1731.5714
$ws.Range("L34").Value = 3552.75
$ws.Range("M34").Value = -1529.5714
$ws.Range("N34").Value = -3956.75
$ws.Range("H99").Value = 2873.25
$ws.Range("I99").Value = 1936.3334
$ws.Range("J99").Value = 3435.4
$ws.Range("K99").Value = 1936.3334
$ws.Range("L99").Value = 3435.4
$ws.Range("M99").Value = -438.3334
$ws.Range("N99").Value = -6431.4
$ws.Range("H122").Value = 3016.1428
$ws.Range("I122").Value = 2999.5
$ws.Range("K122").Value = 8998.5
$ws.Range("M122").Value = -6548.5
$ws.Range("H126").Value = 2873.25
$ws.Range("I126").Value = 1936.3334
$ws.Range("J126").Value = 3435.4
$ws.Range("K126").Value = 5809.0002
$ws.Range("L126").Value = 10306.2
$ws.Range("M126").Value = -3339.0002
$ws.Range("N126").Value = -15246.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 352.57144
$ws.Range("I97").Value = 223.66667
$ws.Range("J97").Value = 449.25
$ws.Range("K97").Value = 671.00001
$ws.Range("L97").Value = 1347.75
$ws.Range("M97").Value = -175.00001
$ws.Range("N97").Value = -2339.75
$ws.Range("H109").Value = 1906.25
$ws.Range("I109").Value = 1750
$ws.Range("K109").Value = 5250
$ws.Range("M109").Value = -4210

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 500
$ws.Range("I3").Value = 500
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 500
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -384
$ws.Range("N3").ClearContents()
$ws.Range("H98").Value = 11549.8
$ws.Range("J98").Value = 11549.8
$ws.Range("L98").Value = 11549.8
$ws.Range("N98").Value = -17539.8
$ws.Range("H122").Value = 2213
$ws.Range("J122").Value = 4999.5
$ws.Range("L122").Value = 14998.5
$ws.Range("N122").Value = -19898.5
$ws.Range("H132").Value = 28852.055
$ws.Range("I132").Value = 37264.57
$ws.Range("J132").Value = 2679.7778
$ws.Range("K132").Value = 111793.71
$ws.Range("L132").Value = 8039.3334
$ws.Range("M132").Value = -109263.71
$ws.Range("N132").Value = -13099.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3469.762
$ws.Range("J40").Value = 5874.75
$ws.Range("L40").Value = 5874.75
$ws.Range("N40").Value = -6146.75
$ws.Range("H55").Value = 1096.3478
$ws.Range("I55").Value = 778.0625
$ws.Range("K55").Value = 778.0625
$ws.Range("M55").Value = -605.0625
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H140").Value = 89000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 89000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 89000
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -99360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 18963
$ws.Range("I52").Value = 18963
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 18963
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -18737
$ws.Range("N52").ClearContents()
$ws.Range("H62").Value = 111451.55
$ws.Range("I62").Value = 4500
$ws.Range("J62").Value = 122146.7
$ws.Range("K62").Value = 4500
$ws.Range("L62").Value = 122146.7
$ws.Range("M62").Value = -3876
$ws.Range("N62").Value = -123394.7
$ws.Range("H65").Value = 111451.55
$ws.Range("I65").Value = 4500
$ws.Range("J65").Value = 122146.7
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 610733.5
$ws.Range("M65").Value = -19380
$ws.Range("N65").Value = -616973.5
$ws.Range("H96").Value = 1433.6
$ws.Range("I96").Value = 1433.6
$ws.Range("K96").Value = 1433.6
$ws.Range("M96").Value = -60.59999999999991
$ws.Range("H99").Value = 74738
$ws.Range("J99").Value = 99476
$ws.Range("L99").Value = 99476
$ws.Range("N99").Value = -105466
$ws.Range("H132").Value = 44914.34
$ws.Range("I132").Value = 50393.617
$ws.Range("K132").Value = 151180.851
$ws.Range("M132").Value = -148650.851
